$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: reword the internal-functionality description and split the
# remainder of the sentence into a new C4 cell.
$ws.Range("B4").Value = "Communicate with the pop up score to tell it what value to display when  the player collides with"

# Row 2: trim the trailing " etc." from the role description.
$ws.Range("B2").Value = "The role of the game manager is to keep track of score values, death parmeters, level diffculty"

# Row 4 (continued): the rest of the new sentence goes into C4.
$ws.Range("C4").Value = "an item. Communicate with the player when it collides with an enemy to display the death."

# Row 2 (continued): a new sentence about instantiating game objects goes into C2.
$ws.Range("C2").Value = "It also instansiates the game objects and asks the world where to spawn them."

# Row 5: the old continuation text is removed entirely.
$ws.Range("B5").Value = ""

# Row 9: the old "Screen display..." text is removed entirely.
$ws.Range("B9").Value = ""

# Update the active selection to match the saved view state.
$ws.Range("B9").Select()
